$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I ("roboticS1Prep") rows 2-41: replace the text "No" with a real
# boolean FALSE value, formatted to still display as TRUE/FALSE text.
$range = $ws.Range("I2:I41")
$range.Value = $false
$range.NumberFormat = '"TRUE";"TRUE";"FALSE"'

# Update the view/selection: the sheet is now scrolled/selected on column I
# (I2:I41, active cell I2) instead of the old H6:H41 selection.
$ws.Range("I2:I41").Select()

$aw = $excel.ActiveWindow
$aw.DisplayGridlines = $true
$aw.DisplayHeadings = $true
$aw.DisplayZeros = $true
$aw.DisplayRightToLeft = $false
$aw.DisplayOutline = $true
